$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps storing plain text values (e.g. "0.9971")
# instead of being auto-converted to numbers by Excel when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2: Bitcoin
$ws.Range("D2").Value = "26.635.14"
$ws.Range("E2").Value = "  +3.33%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.872.90"
$ws.Range("E3").Value = "  +3.08%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "0.9971"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5: BNB
$ws.Range("D5").Value = "282.01"
$ws.Range("E5").Value = "  +1.19%  "

# Row 6: USDC
$ws.Range("D6").Value = "0.9966"
$ws.Range("E6").Value = "  -0.36%  "

# Row 7: XRP
$ws.Range("D7").Value = "0.5146"
$ws.Range("E7").Value = "  +2.93%  "

# Row 8: Cardano
$ws.Range("D8").Value = "0.3543"
$ws.Range("E8").Value = "  +0.56%  "

# Row 9: OKB
$ws.Range("D9").Value = "45.24"
$ws.Range("E9").Value = "  +2.02%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "0.06891"
$ws.Range("E10").Value = "  +3.64%  "

# Row 11: Solana
$ws.Range("D11").Value = "20.18"
$ws.Range("E11").Value = "  -0.20%  "

# Row 12: Polygon
$ws.Range("D12").Value = "0.8172"
$ws.Range("E12").Value = "  -3.80%  "

# Row 13: TRON
$ws.Range("D13").Value = "0.07764"
$ws.Range("E13").Value = "  -0.69%  "

# Row 14: WrappedEther
$ws.Range("D14").Value = "1.864.31"
$ws.Range("E14").Value = "  +2.59%  "

# Row 15: Litecoin
$ws.Range("D15").Value = "89.62"
$ws.Range("E15").Value = "  +2.10%  "

# Row 16: Polkadot (only E changes)
$ws.Range("E16").Value = "  +1.57%  "

# Row 17: BinanceUSD
$ws.Range("D17").Value = "0.9968"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18: Avalanche
$ws.Range("D18").Value = "14.39"
$ws.Range("E18").Value = "  +3.08%  "

# Row 19: ShibaInu
$ws.Range("D19").Value = "0.000008137"
$ws.Range("E19").Value = "  +1.33%  "

# Row 20: Dai
$ws.Range("D20").Value = "0.9965"
$ws.Range("E20").Value = "  -0.54%  "

# Row 21: WrappedBTC
$ws.Range("D21").Value = "26.668.73"
$ws.Range("E21").Value = "  +3.19%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "4.809"
$ws.Range("E22").Value = "  +1.19%  "

# Row 23: Cosmos (only E changes)
$ws.Range("E23").Value = "  +0.95%  "

# Row 24: Chainlink
$ws.Range("D24").Value = "6.240"
$ws.Range("E24").Value = "  +2.17%  "

# Row 25: LidoDAOToken
$ws.Range("D25").Value = "2.395"
$ws.Range("E25").Value = "  +12.37%  "

# Row 26: Monero
$ws.Range("D26").Value = "144.78"
$ws.Range("E26").Value = "  +2.17%  "

# Row 27 and 28 swap: EthereumClassic <-> Toncoin
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "1.666"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "17.35"
$ws.Range("E28").Value = "  +2.41%  "

# Row 29: BitcoinCash
$ws.Range("D29").Value = "110.79"
$ws.Range("E29").Value = "  +1.62%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("D30").Value = "4.387"
$ws.Range("E30").Value = "  +1.37%  "

# Row 31: Filecoin
$ws.Range("D31").Value = "4.344"
$ws.Range("E31").Value = "  +2.55%  "

# Row 32: Stellar
$ws.Range("D32").Value = "0.08813"
$ws.Range("E32").Value = "  +0.27%  "

# Row 33: Hedera
$ws.Range("D33").Value = "0.04911"
$ws.Range("E33").Value = "  +1.93%  "

# Row 34: ARBITRUM
$ws.Range("D34").Value = "1.176"
$ws.Range("E34").Value = "  +3.75%  "

# Row 35: ImmutableX
$ws.Range("D35").Value = "0.7464"
$ws.Range("E35").Value = "  +0.59%  "

# Row 36: HuobiToken
$ws.Range("D36").Value = "2.867"
$ws.Range("E36").Value = "  -1.60%  "

# Row 37: MXToken
$ws.Range("D37").Value = "3.277"
$ws.Range("E37").Value = "  +6.61%  "

# Row 38: RenderToken
$ws.Range("D38").Value = "2.410"
$ws.Range("E38").Value = "  -3.71%  "

# Row 39: VeChain
$ws.Range("D39").Value = "0.01877"
$ws.Range("E39").Value = "  +0.84%  "

# Row 40: TheSandbox
$ws.Range("D40").Value = "0.5241"
$ws.Range("E40").Value = "  -1.87%  "

# Row 41: TrustWalletToken
$ws.Range("D41").Value = "0.9699"
$ws.Range("E41").Value = "  -0.20%  "

# Row 42: Quant
$ws.Range("D42").Value = "116.75"
$ws.Range("E42").Value = "  +3.25%  "

# Row 43: FraxShare
$ws.Range("D43").Value = "6.297"
$ws.Range("E43").Value = "  +1.08%  "

# Row 44: Aptos
$ws.Range("D44").Value = "8.152"
$ws.Range("E44").Value = "  -0.48%  "

# Row 45: PaxDollar
$ws.Range("D45").Value = "0.9955"
$ws.Range("E45").Value = "  -0.47%  "

# Row 46: Decentraland
$ws.Range("D46").Value = "0.4592"
$ws.Range("E46").Value = "  -2.97%  "

# Row 47: Algorand
$ws.Range("D47").Value = "0.1364"
$ws.Range("E47").Value = "  -1.86%  "

# Row 48: EnergySwap
$ws.Range("D48").Value = "9.508"
$ws.Range("E48").Value = "  +2.80%  "

# Row 49: Elrond
$ws.Range("D49").Value = "36.47"
$ws.Range("E49").Value = "  +1.76%  "

# Row 50: NEARProtocol
$ws.Range("D50").Value = "1.515"
$ws.Range("E50").Value = "  +1.71%  "

# Row 51: Cronos
$ws.Range("D51").Value = "0.05923"
$ws.Range("E51").Value = "  +0.58%  "
